$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4249093333333334
$ws.Range("H2").Value = 1.274728
$ws.Range("I2").Value = 0.06780552842016908
$ws.Range("J2").Value = 0.06780552842016908
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 47.01311929447466
$ws.Range("R2").Value = 423.118073650272
$ws.Range("S2").Value = 0.03713314790860505
$ws.Range("T2").Value = 0.03713314790860505

# Row 3
$ws.Range("G3").Value = 0.4249093333333334
$ws.Range("H3").Value = 1.274728
$ws.Range("I3").Value = 0.06780552842016908
$ws.Range("J3").Value = 0.06780552842016908
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 27.05084496336445
$ws.Range("R3").Value = 243.45760467028
$ws.Range("S3").Value = 0.02136601531980054
$ws.Range("T3").Value = 0.02136601531980054

# Row 4
$ws.Range("G4").Value = 0.4249093333333334
$ws.Range("H4").Value = 1.274728
$ws.Range("I4").Value = 0.06780552842016908
$ws.Range("J4").Value = 0.06780552842016908
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 11.78249843065244
$ws.Range("R4").Value = 106.042485875872
$ws.Range("S4").Value = 0.0093063651917635
$ws.Range("T4").Value = 0.0093063651917635

# Row 5
$ws.Range("I5").Value = 0.4624930683973976
$ws.Range("J5").Value = 0.4624930683973975
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 320.6706341509293
$ws.Range("R5").Value = 2886.035707358364
$ws.Range("S5").Value = 0.2532805792631611
$ws.Range("T5").Value = 0.253280579263161

# Row 6
$ws.Range("I6").Value = 0.4624930683973976
$ws.Range("J6").Value = 0.4624930683973975
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("S6").Value = 0.1457349306895308
$ws.Range("T6").Value = 0.1457349306895308

# Row 7
$ws.Range("I7").Value = 0.4624930683973976
$ws.Range("J7").Value = 0.4624930683973975
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 80.36695501895154
$ws.Range("R7").Value = 723.3025951705639
$ws.Range("S7").Value = 0.06347755844470571
$ws.Range("T7").Value = 0.06347755844470569

# Row 8
$ws.Range("G8").Value = 2.943425333333333
$ws.Range("H8").Value = 8.830276
$ws.Range("I8").Value = 0.4697014031824334
$ws.Range("J8").Value = 0.4697014031824334
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 325.6685496758026
$ws.Range("R8").Value = 2931.016947082223
$ws.Range("S8").Value = 0.2572281653668903
$ws.Range("T8").Value = 0.2572281653668903

# Row 9
$ws.Range("G9").Value = 2.943425333333333
$ws.Range("H9").Value = 8.830276
$ws.Range("I9").Value = 0.4697014031824334
$ws.Range("J9").Value = 0.4697014031824334
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 187.3861930229178
$ws.Range("R9").Value = 1686.47573720626
$ws.Range("S9").Value = 0.1480063294240552
$ws.Range("T9").Value = 0.1480063294240552

# Row 10
$ws.Range("G10").Value = 2.943425333333333
$ws.Range("H10").Value = 8.830276
$ws.Range("I10").Value = 0.4697014031824334
$ws.Range("J10").Value = 0.4697014031824334
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.18832399999999
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 81.61954009971377
$ws.Range("R10").Value = 734.5758608974239
$ws.Range("S10").Value = 0.06446690839148794
$ws.Range("T10").Value = 0.06446690839148793
